$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: rewrite a paragraph's body (everything except the trailing
# paragraph mark) as a sequence of "lines" joined by manual line breaks
# (Word's wdLineBreak / OOXML <w:br w:type="textWrapping"/>), then apply the
# "Verbatim Char" character style to each line's text (but not to the break
# characters themselves, so Word keeps the breaks in their own runs).
# Finally the paragraph style is switched to "Source Code".
# ---------------------------------------------------------------------------
function Set-CodeBlock($paraIndex, [string[]]$lines) {
    $p = $d.Paragraphs.Item($paraIndex)
    $pRange = $p.Range
    $bodyStart = $pRange.Start
    $bodyEnd = $pRange.End - 1   # exclude the paragraph mark

    $body = $d.Range($bodyStart, $bodyEnd)
    $nl = [char]11
    $joined = [string]::Join($nl, $lines)
    $body.Text = $joined

    $pos = $bodyStart
    foreach ($line in $lines) {
        $segStart = $pos
        $segEnd = $pos + $line.Length
        $seg = $d.Range($segStart, $segEnd)
        $seg.Style = "Verbatim Char"
        $pos = $segEnd + 1
    }

    $p.Style = "Source Code"
}

$topLines = @(
    "TOP VALUES:",
    "law 37.85900298929792",
    "business 32.87730435677244",
    "enforcement 27.661956717563825",
    "there 26.71258693822815",
    "no 21.023840498701418",
    "we 20.368952478920193",
    "freedom 19.949151468604846",
    "will 19.711108571025164",
    "make 19.27829134912094",
    "do 18.29887958039273",
    "negro 17.7213002579283",
    "amendment 17.662006465968325",
    "accept 17.01798088257794",
    "congress 16.84090924743422",
    "america 16.585253111580684",
    "method 14.88972168638792",
    "islands 14.88972168638792",
    "arbitration 14.88972168638792",
    "south 14.304357332097728",
    "prayer 13.7751544999991"
)

$bottomLines = @(
    "BOTTOM VALUES:",
    "which -25.587694865809624",
    "been -24.768251604369084",
    "powerful -24.257396715885317",
    "powers -20.639945618320855",
    "union -19.808186401517098",
    "myself -19.68520084932692",
    "me -19.212377774373138",
    "foreign -18.94148581831667",
    "my -18.87101149359202",
    " -18.74433136968643",
    "on -18.413088656786343",
    "opinion -16.987689023401476",
    "happy -16.890688937556916",
    "fellow-citizens -16.81270986608169",
    "spirit -16.686359964373743",
    "period -16.515374263165533",
    "limits -16.276798224583224",
    "measures -15.846495944172565",
    "country's -15.213400850163701",
    "democracy -14.591557939972915"
)

# Find the "TOP VALUES:" / "BOTTOM VALUES:" paragraphs and the two
# paragraphs that follow each of them (which get re-styled to
# "First Paragraph"). Locate them by content so this is resilient to
# paragraph-index drift.
$topIndex = 0
$bottomIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("TOP VALUES:")) { $topIndex = $i }
    if ($t.StartsWith("BOTTOM VALUES:")) { $bottomIndex = $i }
}

$afterTopIndex = $topIndex + 1
$afterBottomIndex = $bottomIndex + 1

# Apply the code-block formatting (do the later paragraph first so the
# earlier paragraph's index/position is unaffected by this one's edits).
Set-CodeBlock $bottomIndex $bottomLines
Set-CodeBlock $topIndex $topLines

# The paragraphs immediately following each code block switch from
# "Body Text" to "First Paragraph".
$d.Paragraphs.Item($afterTopIndex).Style = "First Paragraph"
$d.Paragraphs.Item($afterBottomIndex).Style = "First Paragraph"
